$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "64.304.66"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "3.501.54"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "3.502.90"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.374"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("D13").Value = "4.104.62"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000179"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "3.510.71"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.16%  "
$ws.Range("D18").Value = "64.313.74"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.569"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "3.642.11"
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  +3.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").Value = "3.522.42"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.147"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "23.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "164.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0782"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.808"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").Value = "2.476.27"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.917"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.54%  "
